$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell updates. NumberFormat is temporarily forced to text ("@")
# before assigning values so Excel does not auto-convert numeric-looking
# strings (e.g. "1.00", "0.593") into actual numbers, then the original
# style is restored so cell formatting is unchanged.
function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "69.283.16"
Set-TextValue $ws.Range("E2") "  +1.62%  "
Set-TextValue $ws.Range("D3") "3.396.60"
Set-TextValue $ws.Range("E3") "  +1.47%  "
Set-TextValue $ws.Range("D4") "1.00"
Set-TextValue $ws.Range("E4") "  +0.07%  "
Set-TextValue $ws.Range("D5") "580.92"
Set-TextValue $ws.Range("E5") "  -0.37%  "
Set-TextValue $ws.Range("D6") "179.02"
Set-TextValue $ws.Range("E6") "  +1.22%  "
Set-TextValue $ws.Range("E7") "  -0.02%  "
Set-TextValue $ws.Range("D8") "0.593"
Set-TextValue $ws.Range("E8") "  +0.45%  "
Set-TextValue $ws.Range("D9") "0.198"
Set-TextValue $ws.Range("E9") "  +8.61%  "
Set-TextValue $ws.Range("D10") "0.586"
Set-TextValue $ws.Range("E10") "  +0.97%  "
Set-TextValue $ws.Range("D11") "48.30"
Set-TextValue $ws.Range("E11") "  +0.67%  "
Set-TextValue $ws.Range("D12") "0.0000283"
Set-TextValue $ws.Range("E12") "  +3.71%  "
Set-TextValue $ws.Range("D13") "682.72"
Set-TextValue $ws.Range("E13") "  -0.37%  "
Set-TextValue $ws.Range("D14") "8.59"
Set-TextValue $ws.Range("E14") "  +2.12%  "
Set-TextValue $ws.Range("D15") "3.937.78"
Set-TextValue $ws.Range("E15") "  +1.35%  "
Set-TextValue $ws.Range("D16") "69.340.42"
Set-TextValue $ws.Range("E16") "  +1.71%  "
Set-TextValue $ws.Range("E17") "  +0.67%  "
Set-TextValue $ws.Range("D18") "3.383.31"
Set-TextValue $ws.Range("E18") "  +1.19%  "
Set-TextValue $ws.Range("D19") "17.68"
Set-TextValue $ws.Range("E19") "  +1.42%  "
Set-TextValue $ws.Range("D20") "11.29"
Set-TextValue $ws.Range("E20") "  +0.88%  "
Set-TextValue $ws.Range("D21") "0.909"
Set-TextValue $ws.Range("E21") "  +1.82%  "
Set-TextValue $ws.Range("D22") "5.36"
Set-TextValue $ws.Range("E22") "  -1.94%  "
Set-TextValue $ws.Range("D23") "17.05"
Set-TextValue $ws.Range("E23") "  +0.65%  "
Set-TextValue $ws.Range("D24") "101.17"
Set-TextValue $ws.Range("E24") "  +0.79%  "
Set-TextValue $ws.Range("D25") "3.89"
Set-TextValue $ws.Range("E25") "  -0.43%  "
Set-TextValue $ws.Range("D26") "2.70"
Set-TextValue $ws.Range("E26") "  +0.34%  "
Set-TextValue $ws.Range("D27") "9.73"
Set-TextValue $ws.Range("E27") "  +2.44%  "
Set-TextValue $ws.Range("D28") "33.49"
Set-TextValue $ws.Range("E28") "  +1.51%  "
Set-TextValue $ws.Range("D29") "8.75"
Set-TextValue $ws.Range("E29") "  +2.84%  "
Set-TextValue $ws.Range("E30") "  +0.06%  "
Set-TextValue $ws.Range("D31") "3.80"
Set-TextValue $ws.Range("E31") "  +15.09%  "
Set-TextValue $ws.Range("B32") "Bittensor"
Set-TextValue $ws.Range("C32") "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws.Range("D32") "554.93"
Set-TextValue $ws.Range("E32") "  -1.22%  "
Set-TextValue $ws.Range("B33") "Cosmos"
Set-TextValue $ws.Range("C33") "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D33") "11.03"
Set-TextValue $ws.Range("E33") "  -0.34%  "
Set-TextValue $ws.Range("E34") "  +0.10%  "
Set-TextValue $ws.Range("D35") "57.98"
Set-TextValue $ws.Range("E35") "  -0.02%  "
Set-TextValue $ws.Range("E36") "  +0.09%  "
Set-TextValue $ws.Range("D37") "3.607.33"
Set-TextValue $ws.Range("E37") "  -2.95%  "
Set-TextValue $ws.Range("D38") "0.141"
Set-TextValue $ws.Range("E38") "  +2.65%  "
Set-TextValue $ws.Range("E39") "  +1.75%  "
Set-TextValue $ws.Range("D40") "0.0₃0746"
Set-TextValue $ws.Range("E40") "  +11.24%  "
Set-TextValue $ws.Range("D41") "3.30"
Set-TextValue $ws.Range("E41") "  +4.66%  "
Set-TextValue $ws.Range("D42") "2.70"
Set-TextValue $ws.Range("E42") "  +3.67%  "
Set-TextValue $ws.Range("E43") "  +5.04%  "
Set-TextValue $ws.Range("D44") "0.0425"
Set-TextValue $ws.Range("E44") "  +3.55%  "
Set-TextValue $ws.Range("E45") "  +0.05%  "
Set-TextValue $ws.Range("E46") "  +1.27%  "
Set-TextValue $ws.Range("D47") "0.129"
Set-TextValue $ws.Range("E47") "  +0.25%  "
Set-TextValue $ws.Range("D50") "130.92"
Set-TextValue $ws.Range("E50") "  -0.33%  "
Set-TextValue $ws.Range("D51") "2.61"
Set-TextValue $ws.Range("E51") "  +1.95%  "
